# Split three paragraphs' single runs into one run per word (and one run
# per inter-word space), while leaving the paragraph's own text unchanged.
#
# We can't just InsertAfter/collapse+retype piece by piece, because this
# host re-merges adjacent runs that end up with identical formatting.
# Instead we target each paragraph's whole text with Range.InsertXML and
# feed it a tiny WordprocessingML package whose body is the replacement
# run sequence - no <w:pPr> is included, so the paragraph keeps its
# existing style, and each <w:r> stays a distinct run with no empty
# <w:rPr/> left behind.

$d = $word.ActiveDocument

function Split-IntoWordRuns([string]$text) {
    # Break on spaces, keeping each space as its own token, mirroring the
    # target markup (one run per word, one run per single space).
    $tokens = New-Object System.Collections.ArrayList
    $cur = ""
    foreach ($ch in $text.ToCharArray()) {
        if ($ch -eq ' ') {
            if ($cur.Length -gt 0) { [void]$tokens.Add($cur); $cur = "" }
            [void]$tokens.Add(" ")
        } else {
            $cur += $ch
        }
    }
    if ($cur.Length -gt 0) { [void]$tokens.Add($cur) }
    return $tokens
}

function Build-RunsXml([string]$text) {
    $tokens = Split-IntoWordRuns $text
    $sb = New-Object System.Text.StringBuilder
    foreach ($tok in $tokens) {
        [void]$sb.Append('<w:r><w:t xml:space="preserve">')
        [void]$sb.Append($tok)
        [void]$sb.Append('</w:t></w:r>')
    }
    return $sb.ToString()
}

function Set-ParagraphWordRuns([string]$oldText) {
    $found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false,
                                      $false, $true, 1, $false, "", 0)
    $rng = $d.Content
    $rng.Start = $d.Content.Start
    $start = $d.Content.Find.Parent.Start
    # Locate the exact range again via a fresh Find so we have solid
    # Start/End offsets to build a Range from.
    $search = $d.Range(0, $d.Content.End)
    $search.Find.Execute($oldText, $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0) | Out-Null
    $target = $d.Range($search.Start, $search.End)

    $runsXml = Build-RunsXml $oldText
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" ' +
           'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $runsXml + '</w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xml)
}

Set-ParagraphWordRuns "Questions: Vector addition and scalar multiplication"
Set-ParagraphWordRuns "Renee Knapp, Kin Wang Pang"
Set-ParagraphWordRuns "A selection of questions for the study guide on vector addition and scalar multiplication."
